# Updates cryptos list (Price / Volume(1h) columns) per the
# "Updated cryptos list" GitHub Actions automation commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Rows whose new price string is unambiguously numeric would be
# auto-coerced to a Number by Excel's normal type inference, so for
# those we briefly force text format, assign, then restore the
# cell's style to Normal (keeps the General/default numFmt+no style
# index while preserving the Text cell type).
$ws.Range("D2").Value = '63.971.75'
$ws.Range("D3").Value = '3.322.34'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = '3.316.74'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.619'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '52.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '3.858.19'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = '3.324.26'
$ws.Range("D20").Value = '63.886.31'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.969'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '424.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.59'
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '593.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("D42").Value = '3.086.13'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.13'
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume(1h)) updates ---
$ws.Range("E2").Value = '  -3.48%  '
$ws.Range("E3").Value = '  -5.63%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -4.29%  '
$ws.Range("E6").Value = '  -4.13%  '
$ws.Range("E7").Value = '  -3.93%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -5.65%  '
$ws.Range("E10").Value = '  -2.91%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("E12").Value = '  -5.36%  '
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("E14").Value = '  -3.89%  '
$ws.Range("E15").Value = '  -5.45%  '
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("E17").Value = '  -3.97%  '
$ws.Range("E18").Value = '  -5.39%  '
$ws.Range("E19").Value = '  -3.04%  '
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("E21").Value = '  -3.96%  '
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("E23").Value = '  +7.76%  '
$ws.Range("E24").Value = '  -3.17%  '
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("E28").Value = '  -2.51%  '
$ws.Range("E29").Value = '  -6.13%  '
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("E32").Value = '  -8.66%  '
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -4.39%  '
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  -9.18%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -5.83%  '
$ws.Range("E40").Value = '  -7.25%  '
$ws.Range("E41").Value = '  -5.17%  '
$ws.Range("E42").Value = '  -5.26%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  -5.28%  '
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("E49").Value = '  -4.92%  '
$ws.Range("E50").Value = '  -2.96%  '
$ws.Range("E51").Value = '  -6.86%  '

